$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename portraitId -> eventAction (npcNameId stays the same)
$ws.Range("E1").Value = "eventAction"

# Row 9 (buildingId=6) and Row 11 (buildingId=8) swap their comment/name values
$ws.Range("B9").Value = "事务所"
$ws.Range("B11").Value = "旅馆"

# eventAction column values recoded to the new lowercase taxonomy
$ws.Range("E3").Value = "government"
$ws.Range("E4").Value = "government"
$ws.Range("E5").Value = "tarven"
$ws.Range("E6").Value = "plaza"
$ws.Range("E7").Value = "exchange"
$ws.Range("E8").Value = "shipyard"
$ws.Range("E9").Value = "shop"
$ws.Range("E10").Value = "dock"
$ws.Range("E11").Value = "inn"
$ws.Range("E12").Value = "relic"
$ws.Range("E13").Value = "relic"
$ws.Range("E14").Value = "relic"

# Selection moves from G2 to D12
$ws.Range("D12").Select()
